# prep for Feb 18
# Updates the cached "Date Placeholder" field text across the slide master
# and every slide layout (6/10/21 -> 2/17/22), and makes two content edits:
#   - Slide 6 content placeholder: "Basic forking" -> "Adding actions"
#     (and drop one of the blank trailing paragraphs)
#   - Slide 15 content placeholder: rewrite the forking blurb into the
#     "add a Github action" sentence, split across three runs.

function Update-DateText {
    param($shapes, $oldText, $newText)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $found = $tr.Find($oldText, 0)
            if ($found) {
                $found.Text = $newText
            }
        }
    }
}

$p = $ppt.ActivePresentation
$oldDate = "6/10/21"
$newDate = "2/17/22"

# --- Update the cached date text on the slide master ---
Update-DateText $p.SlideMaster.Shapes $oldDate $newDate

# --- Update the cached date text on every slide layout ---
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DateText $layout.Shapes $oldDate $newDate
}

# --- Slide 6: "Basic forking" -> "Adding actions", drop one blank paragraph ---
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(2)
$tr6 = $shp6.TextFrame.TextRange
$found6 = $tr6.Find("Basic forking", 0)
$found6.Text = "Adding actions"

$tr6b = $shp6.TextFrame.TextRange
$blankPara = $tr6b.Paragraphs(6, 1)
$blankPara.Delete()

# --- Slide 15: rewrite the sentence into three runs ---
$s15 = $p.Slides.Item(15)
$shp15 = $s15.Shapes.Item(2)
$tr15 = $shp15.TextFrame.TextRange
$tr15.Text = ""
$null = $tr15.InsertAfter("We are going to add a ")
$null = $tr15.InsertAfter("Github")
$null = $tr15.InsertAfter(" action to thank collaborators")
